# Normalize "Recorded By" (column G) values: when the system/System account
# appears first in the comma-separated list of recorders, move it to the end
# (equivalently, rotate the last entry of the list to the front) so that the
# list starts with the non-system recorder(s).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value()

    if ($null -eq $val) { continue }

    $text = [string]$val
    if ($text -eq "") { continue }

    $parts = $text -split ", "

    if ($parts.Count -gt 1 -and ($parts[0] -eq "System" -or $parts[0] -eq "system")) {
        $rotated = @($parts[-1]) + $parts[0..($parts.Count - 2)]
        $newText = $rotated -join ", "
        $cell.Value = $newText
    }
}
